# Updates cryptos list values to reflect latest coinranking.com data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'28.026.17"
$ws.Range('E2').Value = "'  -0.77%  "
$ws.Range('D3').Value = "'1.802.01"
$ws.Range('E3').Value = "'  -0.20%  "
$ws.Range('D4').Value = "'1.002"
$ws.Range('E4').Value = "'  -0.05%  "
$ws.Range('D5').Value = "'310.10"
$ws.Range('E5').Value = "'  -1.47%  "
$ws.Range('E6').Value = "'  -0.04%  "
$ws.Range('D7').Value = "'0.5078"
$ws.Range('E7').Value = "'  -3.34%  "
$ws.Range('D8').Value = "'0.3853"
$ws.Range('E8').Value = "'  +0.47%  "
$ws.Range('D9').Value = "'0.09006"
$ws.Range('E9').Value = "'  +12.45%  "
$ws.Range('D10').Value = "'1.095"
$ws.Range('E10').Value = "'  -0.77%  "
$ws.Range('D11').Value = "'40.80"
$ws.Range('E11').Value = "'  -1.54%  "
$ws.Range('D12').Value = "'6.379"
$ws.Range('E12').Value = "'  +0.33%  "
$ws.Range('D13').Value = "'1.002"
$ws.Range('E13').Value = "'  -0.01%  "
$ws.Range('D14').Value = "'20.34"
$ws.Range('E14').Value = "'  -1.47%  "
$ws.Range('B15').Value = "'WrappedEther"
$ws.Range('C15').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D15').Value = "'1.803.46"
$ws.Range('E15').Value = "'  +0.07%  "
$ws.Range('B16').Value = "'Chainlink"
$ws.Range('C16').Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range('D16').Value = "'7.300"
$ws.Range('E16').Value = "'  -0.91%  "
$ws.Range('D17').Value = "'0.00001115"
$ws.Range('E17').Value = "'  +1.45%  "
$ws.Range('D18').Value = "'92.24"
$ws.Range('E18').Value = "'  -0.55%  "
$ws.Range('D19').Value = "'0.06575"
$ws.Range('E19').Value = "'  -0.53%  "
$ws.Range('E20').Value = "'  -0.08%  "
$ws.Range('D21').Value = "'17.23"
$ws.Range('E21').Value = "'  -1.02%  "
$ws.Range('D22').Value = "'6.009"
$ws.Range('E22').Value = "'  +0.54%  "
$ws.Range('D23').Value = "'28.048.99"
$ws.Range('E23').Value = "'  -0.91%  "
$ws.Range('D24').Value = "'11.06"
$ws.Range('E24').Value = "'  -1.35%  "
$ws.Range('D25').Value = "'2.221"
$ws.Range('E25').Value = "'  -0.63%  "
$ws.Range('D26').Value = "'158.30"
$ws.Range('E26').Value = "'  -1.08%  "
$ws.Range('D27').Value = "'2.005.51"
$ws.Range('E27').Value = "'  -0.25%  "
$ws.Range('D28').Value = "'2.404"
$ws.Range('E28').Value = "'  +1.01%  "
$ws.Range('D29').Value = "'20.33"
$ws.Range('E29').Value = "'  -0.88%  "
$ws.Range('D30').Value = "'127.26"
$ws.Range('E30').Value = "'  +3.31%  "
$ws.Range('D31').Value = "'0.1087"
$ws.Range('E31').Value = "'  +0.24%  "
$ws.Range('D32').Value = "'1.047"
$ws.Range('E32').Value = "'  -1.28%  "
$ws.Range('D33').Value = "'5.574"
$ws.Range('E33').Value = "'  +0.00%  "
$ws.Range('D34').Value = "'3.643"
$ws.Range('E34').Value = "'  -0.86%  "
$ws.Range('D35').Value = "'0.06917"
$ws.Range('E35').Value = "'  -4.88%  "
$ws.Range('D36').Value = "'8.975"
$ws.Range('E36').Value = "'  +0.68%  "
$ws.Range('D37').Value = "'0.02331"
$ws.Range('E37').Value = "'  +0.38%  "
$ws.Range('D38').Value = "'0.2167"
$ws.Range('E38').Value = "'  -0.17%  "
$ws.Range('D39').Value = "'4.982"
$ws.Range('E39').Value = "'  -3.80%  "
$ws.Range('D40').Value = "'11.40"
$ws.Range('E40').Value = "'  -8.53%  "
$ws.Range('D41').Value = "'0.6114"
$ws.Range('E41').Value = "'  -1.84%  "
$ws.Range('E42').Value = "'  +0.02%  "
$ws.Range('D43').Value = "'1.153"
$ws.Range('E43').Value = "'  -1.29%  "
$ws.Range('D44').Value = "'13.22"
$ws.Range('E44').Value = "'  -0.55%  "
$ws.Range('D45').Value = "'1.293"
$ws.Range('E45').Value = "'  -5.82%  "
$ws.Range('D46').Value = "'0.5889"
$ws.Range('E46').Value = "'  -2.51%  "
$ws.Range('D47').Value = "'3.706"
$ws.Range('D48').Value = "'124.34"
$ws.Range('E48').Value = "'  -2.24%  "
$ws.Range('D49').Value = "'1.933"
$ws.Range('E49').Value = "'  -0.17%  "
$ws.Range('D50').Value = "'1.183"
$ws.Range('E50').Value = "'  -2.48%  "
$ws.Range('D51').Value = "'0.06733"
$ws.Range('E51').Value = "'  -1.69%  "
